$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty clock-out time and duration for row 19
$ws.Range("C19").Value = "22:37:50"
$ws.Range("D19").Value = "0.99 Hours"

# Duplicate the formatting of row 19 onto the new total row, then set its values
$ws.Range("C19:D19").Copy()
$ws.Range("C20:D20").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C20").Value = "Total Duration:"
$ws.Range("D20").Value = "26.5 Hours"

$excel.CutCopyMode = 0
